$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("A2").Value = "total"
$ws.Range("B2").Value = 5752578.84
$ws.Range("C2").Value = 5.36
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0.72
$ws.Range("F2").Value = 3.91
$ws.Range("G2").Value = 3.65
$ws.Range("H2").Value = 1.07

# Delete row 3 entirely (shifts rows up, removing the extra row)
$ws.Rows.Item(3).Delete()
